$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 266.75
$ws.Range("J6").Value = 394.75
$ws.Range("L6").Value = 1184.25
$ws.Range("N6").Value = -1408.25

$ws.Range("H43").Value = 1189.2632
$ws.Range("I43").Value = 1008.4
$ws.Range("K43").Value = 1008.4
$ws.Range("M43").Value = -939.4

$ws.Range("H48").Value = 592.5
$ws.Range("I48").Value = 566.6667
$ws.Range("J48").Value = 670
$ws.Range("K48").Value = 1700.0001
$ws.Range("L48").Value = 2010
$ws.Range("M48").Value = -1408.0001
$ws.Range("N48").Value = -2594

$ws.Range("H56").Value = 592.5
$ws.Range("I56").Value = 566.6667
$ws.Range("J56").Value = 670
$ws.Range("K56").Value = 1700.0001
$ws.Range("L56").Value = 2010
$ws.Range("M56").Value = -1166.0001
$ws.Range("N56").Value = -3078

$ws.Range("H132").Value = 4698.75
$ws.Range("I132").Value = 5029.6895
$ws.Range("K132").Value = 15089.0685
$ws.Range("M132").Value = -12559.0685

$ws.Range("H138").Value = 4980.964
$ws.Range("I138").Value = 1559.174
$ws.Range("J138").Value = 7440.375
$ws.Range("K138").Value = 4677.522
$ws.Range("L138").Value = 22321.125
$ws.Range("M138").Value = 462.4780000000001
$ws.Range("N138").Value = -32601.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4494.5
$ws.Range("I32").Value = 3866.149
$ws.Range("K32").Value = 3866.149
$ws.Range("M32").Value = -3579.149

$ws.Range("H74").Value = 26346998
$ws.Range("I74").Value = 31286108
$ws.Range("J74").Value = 5082.1665
$ws.Range("K74").Value = 31286108
$ws.Range("L74").Value = 5082.1665
$ws.Range("M74").Value = -31285234
$ws.Range("N74").Value = -6830.1665

$ws.Range("H77").Value = 26346998
$ws.Range("I77").Value = 31286108
$ws.Range("J77").Value = 5082.1665
$ws.Range("K77").Value = 156430540
$ws.Range("L77").Value = 25410.8325
$ws.Range("M77").Value = -156426172
$ws.Range("N77").Value = -34146.8325

$ws.Range("H122").Value = 37039336
$ws.Range("I122").Value = 2242.5715
$ws.Range("K122").Value = 6727.7145
$ws.Range("M122").Value = -4277.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 44317.43
$ws.Range("I86").Value = 21746.4
$ws.Range("K86").Value = 21746.4
$ws.Range("M86").Value = -20623.4

$ws.Range("H89").Value = 44317.43
$ws.Range("I89").Value = 21746.4
$ws.Range("K89").Value = 108732
$ws.Range("M89").Value = -103116

$ws.Range("H99").Value = 2766.8823
$ws.Range("I99").Value = 2130.0908
$ws.Range("J99").Value = 3934.3333
$ws.Range("K99").Value = 2130.0908
$ws.Range("L99").Value = 3934.3333
$ws.Range("M99").Value = -632.0907999999999
$ws.Range("N99").Value = -6930.3333

$ws.Range("H107").Value = 1587.8572
$ws.Range("I107").Value = 1535
$ws.Range("J107").Value = 1627.5
$ws.Range("K107").Value = 1535
$ws.Range("L107").Value = 1627.5
$ws.Range("M107").Value = 385
$ws.Range("N107").Value = -5467.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6947814.5
$ws.Range("I31").Value = 2467.1072
$ws.Range("J31").Value = 31256530
$ws.Range("K31").Value = 2467.1072
$ws.Range("L31").Value = 31256530
$ws.Range("M31").Value = -2172.1072
$ws.Range("N31").Value = -31257120

$ws.Range("H34").Value = 6947814.5
$ws.Range("I34").Value = 2467.1072
$ws.Range("J34").Value = 31256530
$ws.Range("K34").Value = 2467.1072
$ws.Range("L34").Value = 31256530
$ws.Range("M34").Value = -2265.1072
$ws.Range("N34").Value = -31256934

$ws.Range("H99").Value = 8142.857
$ws.Range("I99").Value = 8666.666999999999
$ws.Range("K99").Value = 8666.666999999999
$ws.Range("M99").Value = -7168.666999999999

$ws.Range("H126").Value = 8142.857
$ws.Range("I126").Value = 8666.666999999999
$ws.Range("K126").Value = 26000.001
$ws.Range("M126").Value = -23530.001

$ws.Range("H132").Value = 66979.09
$ws.Range("I132").Value = 73691.96000000001
$ws.Range("J132").Value = 4325.6665
$ws.Range("K132").Value = 221075.88
$ws.Range("L132").Value = 12976.9995
$ws.Range("M132").Value = -218545.88
$ws.Range("N132").Value = -18036.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2975.5715
$ws.Range("J5").Value = 2975.5715
$ws.Range("L5").Value = 8926.7145
$ws.Range("N5").Value = -9150.7145

$ws.Range("H75").Value = 1962.6
$ws.Range("I75").Value = 1906.5
$ws.Range("K75").Value = 5719.5
$ws.Range("M75").Value = -4721.5

$ws.Range("H78").Value = 1962.6
$ws.Range("I78").Value = 1906.5
$ws.Range("K78").Value = 17158.5
$ws.Range("M78").Value = -12166.5

$ws.Range("H86").Value = 1061.3077
$ws.Range("I86").Value = 734.2222
$ws.Range("J86").Value = 1797.25
$ws.Range("K86").Value = 2202.6666
$ws.Range("L86").Value = 5391.75
$ws.Range("M86").Value = -1016.6666
$ws.Range("N86").Value = -7763.75

$ws.Range("H89").Value = 1061.3077
$ws.Range("I89").Value = 734.2222
$ws.Range("J89").Value = 1797.25
$ws.Range("K89").Value = 6607.999800000001
$ws.Range("L89").Value = 16175.25
$ws.Range("M89").Value = -679.9998000000005
$ws.Range("N89").Value = -28031.25

$ws.Range("H98").Value = 568.7
$ws.Range("J98").Value = 571.7143
$ws.Range("L98").Value = 1715.1429
$ws.Range("N98").Value = -4711.1429

$ws.Range("H135").Value = 2975.5715
$ws.Range("J135").Value = 2975.5715
$ws.Range("L135").Value = 26780.1435
$ws.Range("N135").Value = -31850.1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12699.2
$ws.Range("I80").Value = 12099.7
$ws.Range("J80").Value = 13898.2
$ws.Range("K80").Value = 12099.7
$ws.Range("L80").Value = 13898.2
$ws.Range("M80").Value = -11101.7
$ws.Range("N80").Value = -15894.2

$ws.Range("H83").Value = 12699.2
$ws.Range("I83").Value = 12099.7
$ws.Range("J83").Value = 13898.2
$ws.Range("K83").Value = 60498.5
$ws.Range("L83").Value = 69491
$ws.Range("M83").Value = -55506.5
$ws.Range("N83").Value = -79475

$ws.Range("H107").Value = 997
$ws.Range("I107").Value = 997
$ws.Range("K107").Value = 997
$ws.Range("M107").Value = 923

$ws.Range("H122").Value = 62501780
$ws.Range("I122").Value = 1125.75
$ws.Range("K122").Value = 3377.25
$ws.Range("M122").Value = -927.25

$ws.Range("H132").Value = 3843.4167
$ws.Range("I132").Value = 3568.6667
$ws.Range("J132").Value = 4667.6665
$ws.Range("K132").Value = 10706.0001
$ws.Range("L132").Value = 14002.9995
$ws.Range("M132").Value = -8176.000100000001
$ws.Range("N132").Value = -19062.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I122").Value = 3105.3555
$ws.Range("K122").Value = 9316.066500000001
$ws.Range("M122").Value = -6866.066500000001

$ws.Range("H136").Value = 6003.0454
$ws.Range("I136").Value = 3471.9285
$ws.Range("K136").Value = 10415.7855
$ws.Range("M136").Value = -7865.7855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1433.2273
$ws.Range("I132").Value = 1164.9459
$ws.Range("J132").Value = 2851.2856
$ws.Range("K132").Value = 3494.8377
$ws.Range("L132").Value = 8553.856800000001
$ws.Range("M132").Value = -964.8377
$ws.Range("N132").Value = -13613.8568
